# Generate Report for Handoff
#
# The "2caa79a4-97ad-499e-a178-6fe2978ed58a.md" file has completed translation
# and is now ready to be handed off, while "4fc6bd29-9a16-425a-97f7-6852e04ca6f4.md"
# remains in translation. The status rows are re-sorted (the newly-ready file
# moves to the bottom row of each table) and its status / handoff file /
# handoff datetime are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "4fc6bd29-9a16-425a-97f7-6852e04ca6f4.md"
$overview.Range("B2").Value = "In Translation"
$overview.Range("C2").Value = "In Translation"
$overview.Range("D2").Value = "2016-18-18 22:18:00"

$overview.Range("A3").Value = "2caa79a4-97ad-499e-a178-6fe2978ed58a.md"
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-18-18 22:18:37"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "4fc6bd29-9a16-425a-97f7-6852e04ca6f4.md"
$zhcn.Range("B2").Value = ".md"
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("D2").Value = "4fc6bd29-9a16-425a-97f7-6852e04ca6f4.d3fe1639c8ab2a2074d3a9e0f6262b790e9425c7.zh-cn.xlf"
$zhcn.Range("E2").Value = "2016-03-18 22:17:58"
$zhcn.Range("H2").Value = "0001-01-01 00:00:00"
$zhcn.Range("I2").Value = "Include"

$zhcn.Range("A3").Value = "2caa79a4-97ad-499e-a178-6fe2978ed58a.md"
$zhcn.Range("B3").Value = ".md"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2caa79a4-97ad-499e-a178-6fe2978ed58a.d7ab2f752541d7b8e5dcf93bd932c789bc177340.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-18 22:18:34"
$zhcn.Range("H3").Value = "0001-01-01 00:00:00"
$zhcn.Range("I3").Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "4fc6bd29-9a16-425a-97f7-6852e04ca6f4.md"
$dede.Range("B2").Value = ".md"
$dede.Range("C2").Value = "In Translation"
$dede.Range("D2").Value = "4fc6bd29-9a16-425a-97f7-6852e04ca6f4.d3fe1639c8ab2a2074d3a9e0f6262b790e9425c7.de-de.xlf"
$dede.Range("E2").Value = "2016-03-18 22:18:00"
$dede.Range("H2").Value = "0001-01-01 00:00:00"
$dede.Range("I2").Value = "Include"

$dede.Range("A3").Value = "2caa79a4-97ad-499e-a178-6fe2978ed58a.md"
$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2caa79a4-97ad-499e-a178-6fe2978ed58a.d7ab2f752541d7b8e5dcf93bd932c789bc177340.de-de.xlf"
$dede.Range("E3").Value = "2016-03-18 22:18:37"
$dede.Range("H3").Value = "0001-01-01 00:00:00"
$dede.Range("I3").Value = "Include"
